$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 35.42516366666666
$ws.Range("H2").Value = 106.275491
$ws.Range("I2").Value = 0.00832770193000585
$ws.Range("J2").Value = 0.008327701930005852
$ws.Range("M2").Value = 0.4394373333333333
$ws.Range("N2").Value = 1.318312
$ws.Range("O2").Value = 0.2944741752765458
$ws.Range("P2").Value = 0.2944741752765458
$ws.Range("Q2").Value = 15.56713945457689
$ws.Range("R2").Value = 140.104255091192
$ws.Range("S2").Value = 0.002452293157787372
$ws.Range("T2").Value = 0.002452293157787372

# Row 3
$ws.Range("G3").Value = 35.42516366666666
$ws.Range("H3").Value = 106.275491
$ws.Range("I3").Value = 0.00832770193000585
$ws.Range("J3").Value = 0.008327701930005852
$ws.Range("O3").Value = 0.4358046333636673
$ws.Range("P3").Value = 0.4358046333636673
$ws.Range("Q3").Value = 23.03845998091944
$ws.Range("R3").Value = 207.346139828275
$ws.Range("S3").Value = 0.003629251086368104
$ws.Range("T3").Value = 0.003629251086368105

# Row 4
$ws.Range("G4").Value = 35.42516366666666
$ws.Range("H4").Value = 106.275491
$ws.Range("I4").Value = 0.00832770193000585
$ws.Range("J4").Value = 0.008327701930005852
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.1607546666666667
$ws.Range("N4").Value = 0.482264
$ws.Range("O4").Value = 0.1077243426939663
$ws.Range("P4").Value = 0.1077243426939663
$ws.Range("Q4").Value = 5.694760376847111
$ws.Range("R4").Value = 51.252843391624
$ws.Range("S4").Value = 0.000897096216561155
$ws.Range("T4").Value = 0.0008970962165611551

# Row 5
$ws.Range("G5").Value = 35.42516366666666
$ws.Range("H5").Value = 106.275491
$ws.Range("I5").Value = 0.00832770193000585
$ws.Range("J5").Value = 0.008327701930005852
$ws.Range("M5").Value = 0.2417443333333333
$ws.Range("N5").Value = 0.725233
$ws.Range("O5").Value = 0.1619968486658205
$ws.Range("P5").Value = 0.1619968486658205
$ws.Range("Q5").Value = 8.563832573822555
$ws.Range("R5").Value = 77.07449316440299
$ws.Range("S5").Value = 0.001349061469289219
$ws.Range("T5").Value = 0.00134906146928922

# Row 6
$ws.Range("I6").Value = 0.01070182047907406
$ws.Range("J6").Value = 0.01070182047907406
$ws.Range("M6").Value = 0.4394373333333333
$ws.Range("N6").Value = 1.318312
$ws.Range("O6").Value = 0.2944741752765458
$ws.Range("P6").Value = 0.2944741752765458
$ws.Range("Q6").Value = 20.00512665028533
$ws.Range("R6").Value = 180.046139852568
$ws.Range("S6").Value = 0.003151409759532981
$ws.Range("T6").Value = 0.003151409759532981

# Row 7
$ws.Range("I7").Value = 0.01070182047907406
$ws.Range("J7").Value = 0.01070182047907406
$ws.Range("O7").Value = 0.4358046333636673
$ws.Range("P7").Value = 0.4358046333636673
$ws.Range("S7").Value = 0.004663902950206655
$ws.Range("T7").Value = 0.004663902950206655

# Row 8
$ws.Range("I8").Value = 0.01070182047907406
$ws.Range("J8").Value = 0.01070182047907406
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.1607546666666667
$ws.Range("N8").Value = 0.482264
$ws.Range("O8").Value = 0.1077243426939663
$ws.Range("P8").Value = 0.1077243426939663
$ws.Range("Q8").Value = 7.318261837010668
$ws.Range("R8").Value = 65.864356533096
$ws.Range("S8").Value = 0.001152846576737081
$ws.Range("T8").Value = 0.00115284657673708

# Row 9
$ws.Range("I9").Value = 0.01070182047907406
$ws.Range("J9").Value = 0.01070182047907406
$ws.Range("M9").Value = 0.2417443333333333
$ws.Range("N9").Value = 0.725233
$ws.Range("O9").Value = 0.1619968486658205
$ws.Range("P9").Value = 0.1619968486658205
$ws.Range("Q9").Value = 11.00526887107633
$ws.Range("R9").Value = 99.04741983968701
$ws.Range("S9").Value = 0.001733661192597339
$ws.Range("T9").Value = 0.001733661192597339

# Row 10
$ws.Range("G10").Value = 51.06824600000001
$ws.Range("H10").Value = 153.204738
$ws.Range("I10").Value = 0.01200505761322374
$ws.Range("J10").Value = 0.01200505761322374
$ws.Range("M10").Value = 0.4394373333333333
$ws.Range("N10").Value = 1.318312
$ws.Range("O10").Value = 0.2944741752765458
$ws.Range("P10").Value = 0.2944741752765458
$ws.Range("Q10").Value = 22.44129384025067
$ws.Range("R10").Value = 201.971644562256
$ws.Range("S10").Value = 0.003535179439801478
$ws.Range("T10").Value = 0.003535179439801479

# Row 11
$ws.Range("G11").Value = 51.06824600000001
$ws.Range("H11").Value = 153.204738
$ws.Range("I11").Value = 0.01200505761322374
$ws.Range("J11").Value = 0.01200505761322374
$ws.Range("O11").Value = 0.4358046333636673
$ws.Range("P11").Value = 0.4358046333636673
$ws.Range("Q11").Value = 33.21180821738334
$ws.Range("R11").Value = 298.9062739564501
$ws.Range("S11").Value = 0.005231859731640674
$ws.Range("T11").Value = 0.005231859731640675

# Row 12
$ws.Range("G12").Value = 51.06824600000001
$ws.Range("H12").Value = 153.204738
$ws.Range("I12").Value = 0.01200505761322374
$ws.Range("J12").Value = 0.01200505761322374
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 0.1607546666666667
$ws.Range("N12").Value = 0.482264
$ws.Range("O12").Value = 0.1077243426939663
$ws.Range("P12").Value = 0.1077243426939663
$ws.Range("Q12").Value = 8.209458862981336
$ws.Range("R12").Value = 73.88512976683201
$ws.Range("S12").Value = 0.001293236940387724
$ws.Range("T12").Value = 0.001293236940387724

# Row 13
$ws.Range("G13").Value = 51.06824600000001
$ws.Range("H13").Value = 153.204738
$ws.Range("I13").Value = 0.01200505761322374
$ws.Range("J13").Value = 0.01200505761322374
$ws.Range("M13").Value = 0.2417443333333333
$ws.Range("N13").Value = 0.725233
$ws.Range("O13").Value = 0.1619968486658205
$ws.Range("P13").Value = 0.1619968486658205
$ws.Range("Q13").Value = 12.34545908377267
$ws.Range("R13").Value = 111.109131753954
$ws.Range("S13").Value = 0.001944781501393863
$ws.Range("T13").Value = 0.001944781501393863

# Row 14
$ws.Range("G14").Value = 4121.876464666667
$ws.Range("H14").Value = 12365.629394
$ws.Range("I14").Value = 0.9689654199776964
$ws.Range("J14").Value = 0.9689654199776964
$ws.Range("M14").Value = 0.4394373333333333
$ws.Range("N14").Value = 1.318312
$ws.Range("O14").Value = 0.2944741752765458
$ws.Range("P14").Value = 0.2944741752765458
$ws.Range("Q14").Value = 1811.306401962547
$ws.Range("R14").Value = 16301.75761766293
$ws.Range("S14").Value = 0.285335292919424
$ws.Range("T14").Value = 0.285335292919424

# Row 15
$ws.Range("G15").Value = 4121.876464666667
$ws.Range("H15").Value = 12365.629394
$ws.Range("I15").Value = 0.9689654199776964
$ws.Range("J15").Value = 0.9689654199776964
$ws.Range("O15").Value = 0.4358046333636673
$ws.Range("P15").Value = 0.4358046333636673
$ws.Range("Q15").Value = 2680.628009825428
$ws.Range("R15").Value = 24125.65208842885
$ws.Range("S15").Value = 0.4222796195954519
$ws.Range("T15").Value = 0.4222796195954519

# Row 16
$ws.Range("G16").Value = 4121.876464666667
$ws.Range("H16").Value = 12365.629394
$ws.Range("I16").Value = 0.9689654199776964
$ws.Range("J16").Value = 0.9689654199776964
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 0.1607546666666667
$ws.Range("N16").Value = 0.482264
$ws.Range("O16").Value = 0.1077243426939663
$ws.Range("P16").Value = 0.1077243426939663
$ws.Range("Q16").Value = 662.6108771186686
$ws.Range("R16").Value = 5963.497894068016
$ws.Range("S16").Value = 0.1043811629602804
$ws.Range("T16").Value = 0.1043811629602804

# Row 17
$ws.Range("G17").Value = 4121.876464666667
$ws.Range("H17").Value = 12365.629394
$ws.Range("I17").Value = 0.9689654199776964
$ws.Range("J17").Value = 0.9689654199776964
$ws.Range("M17").Value = 0.2417443333333333
$ws.Range("N17").Value = 0.725233
$ws.Range("O17").Value = 0.1619968486658205
$ws.Range("P17").Value = 0.1619968486658205
$ws.Range("Q17").Value = 996.4402780332002
$ws.Range("R17").Value = 8967.962502298802
$ws.Range("S17").Value = 0.1569693445025401
$ws.Range("T17").Value = 0.1569693445025401
